# Scheduled runner update: refresh currentAveragePrice / Leve profit figures
# (columns H-N) across the eight crafting-leve worksheets to the latest
# market-board snapshot. No formulas are present in these sheets - every
# figure is a literal numeric value, so each touched cell is written
# directly. Cells that should become empty (e.g. a leve that no longer has
# an HQ profit row) are cleared with ClearContents().

$wb = $excel.ActiveWorkbook

# --- ALC --------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 5 - Met a Sticky End / Animal Glue
$ws.Range("H5").Value = 2309.7778
$ws.Range("I5").Value = 2309.7778
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2309.7778
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2194.7778
$ws.Range("N5").ClearContents()

# Row 11 - Gotta Bounce / Rubber
$ws.Range("H11").Value = 910
$ws.Range("I11").Value = 910
$ws.Range("K11").Value = 910
$ws.Range("M11").Value = -770

# Row 15 - Morning Glass of Ether / Ether
$ws.Range("H15").Value = 369.8
$ws.Range("I15").Value = 369.8
$ws.Range("K15").Value = 1109.4
$ws.Range("M15").Value = -940.4000000000001

# Row 19 - Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 1009.4667
$ws.Range("J19").Value = 1108.8
$ws.Range("L19").Value = 1108.8
$ws.Range("N19").Value = -1458.8

# Row 32 - Automata for the People / Crab Oil
$ws.Range("H32").Value = 985.3125
$ws.Range("I32").Value = 785.375
$ws.Range("J32").Value = 1185.25
$ws.Range("K32").Value = 785.375
$ws.Range("L32").Value = 1185.25
$ws.Range("M32").Value = -459.375
$ws.Range("N32").Value = -1837.25

# Row 44 - Alive and Unwell / Budding Oak Wand
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20924

# Row 63 - Summoning for Dummies / Archaeoskin Codex
$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 20000
$ws.Range("K63").Value = 20000
$ws.Range("M63").Value = -19376

# Row 66 - Summoning the Courage to Be Different (L) / Archaeoskin Codex
$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 20000
$ws.Range("K66").Value = 60000
$ws.Range("M66").Value = -56880

# Row 132 - Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 1705.4736
$ws.Range("I132").Value = 1553.1765
$ws.Range("K132").Value = 4659.529500000001
$ws.Range("M132").Value = -2129.529500000001

# --- ARM ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 63 - Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 2805.875
$ws.Range("I63").Value = 2489.4
$ws.Range("K63").Value = 2489.4
$ws.Range("M63").Value = -1803.4

# Row 66 - A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 2805.875
$ws.Range("I66").Value = 2489.4
$ws.Range("K66").Value = 12447
$ws.Range("M66").Value = -9015

# Row 132 - Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2997
$ws.Range("I132").Value = 2997
$ws.Range("K132").Value = 8991
$ws.Range("M132").Value = -6461

# Row 133 - Shielding My Students / Mountain Chromite Tower Shield
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- BSM ------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 33 - Mors Dagger / Steel Broadsword
$ws.Range("H33").Value = 9642.857
$ws.Range("J33").Value = 10416.667
$ws.Range("L33").Value = 10416.667
$ws.Range("N33").Value = -11088.667

# Row 86 - Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 5279.7
$ws.Range("I86").Value = 4350
$ws.Range("K86").Value = 4350
$ws.Range("M86").Value = -3227

# Row 89 - Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 5279.7
$ws.Range("I89").Value = 4350
$ws.Range("K89").Value = 21750
$ws.Range("M89").Value = -16134

# Row 94 - High Steal / High Steel Nugget
$ws.Range("H94").Value = 2584.8333
$ws.Range("I94").Value = 2584.8333
$ws.Range("K94").Value = 2584.8333
$ws.Range("M94").Value = -2133.8333

# Row 105 - Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 4780.7144
$ws.Range("I105").Value = 4660.3335
$ws.Range("K105").Value = 4660.3335
$ws.Range("M105").Value = -2913.3335

# Row 109 - Here Comes the Hammer / Deepgold Sledgehammer
$ws.Range("H109").Value = 79989.25
$ws.Range("J109").Value = 79989.25
$ws.Range("L109").Value = 79989.25
$ws.Range("N109").Value = -82763.25

# Row 134 - Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 5952.909
$ws.Range("I134").Value = 5248.3
$ws.Range("K134").Value = 15744.9
$ws.Range("M134").Value = -13209.9

# --- CRP ------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 22 - Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150

# Row 99 - O Pine / Pine Lumber
$ws.Range("H99").Value = 4772
$ws.Range("I99").Value = 4112.375
$ws.Range("J99").Value = 5651.5
$ws.Range("K99").Value = 4112.375
$ws.Range("L99").Value = 5651.5
$ws.Range("M99").Value = -2614.375
$ws.Range("N99").Value = -8647.5

# Row 122 - Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 1233.6666
$ws.Range("I122").Value = 1258.2858
$ws.Range("K122").Value = 3774.8574
$ws.Range("M122").Value = -1324.8574

# Row 126 - A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 4772
$ws.Range("I126").Value = 4112.375
$ws.Range("J126").Value = 5651.5
$ws.Range("K126").Value = 12337.125
$ws.Range("L126").Value = 16954.5
$ws.Range("M126").Value = -9867.125
$ws.Range("N126").Value = -21894.5

# Row 130 - Annals of the Empire II / Integral Magitek Rod
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 131 - An Integral Reward / Integral Necklace of Crafting
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

# Row 132 - Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 3999
$ws.Range("I132").Value = 2998
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8994
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6464
$ws.Range("N132").Value = -20060

# --- CUL ------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 39 - Bloody Good Tart, This / Blood Currant Tart
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Row 116 - On a Full Stomach / Sausage Links
$ws.Range("H116").Value = 2000
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# Row 140 - Sweet, Sweet Bean Juice / Mesquite Juice
$ws.Range("H140").Value = 1516.6666
$ws.Range("I140").Value = 1516.6666
$ws.Range("K140").Value = 4549.9998
$ws.Range("M140").Value = 630.0002000000004

# Row 141 - Ocean Explosion / Acqua Pazza
$ws.Range("H141").Value = 2590
$ws.Range("I141").Value = 2590
$ws.Range("K141").Value = 7770
$ws.Range("M141").Value = -2590

# --- GSM ------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 97 - If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 766.3333
$ws.Range("I97").Value = 766.3333
$ws.Range("K97").Value = 766.3333
$ws.Range("M97").Value = -270.3333

# Row 113 - Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 2718.75
$ws.Range("I113").Value = 2025.3334
$ws.Range("J113").Value = 4799
$ws.Range("K113").Value = 2025.3334
$ws.Range("L113").Value = 4799
$ws.Range("M113").Value = 144.6666
$ws.Range("N113").Value = -9139

# --- LTW ------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 46 - Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312

# Row 61 - Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 3111.5
$ws.Range("I61").Value = 3111.5
$ws.Range("K61").Value = 3111.5
$ws.Range("M61").Value = -2909.5

# Row 113 - Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 3111.5
$ws.Range("I113").Value = 3111.5
$ws.Range("K113").Value = 3111.5
$ws.Range("M113").Value = -941.5

# Row 122 - Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3846.25
$ws.Range("I122").Value = 3628.3333
$ws.Range("K122").Value = 10884.9999
$ws.Range("M122").Value = -8434.999899999999

# Row 132 - Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 10811
$ws.Range("J132").Value = 12935.4
$ws.Range("L132").Value = 38806.2
$ws.Range("N132").Value = -43866.2

# --- WVR ------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 3 - Trew Enough / Hempen Chausses
$ws.Range("H3").Value = 1000000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# Row 41 - Half Is the New Double / Linen Halfgloves
$ws.Range("H41").Value = 19855.4
$ws.Range("J41").Value = 19855.4
$ws.Range("L41").Value = 19855.4
$ws.Range("N41").Value = -20635.4

# Row 46 - Crunching the Numbers / Linen Hat
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Row 107 - Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 3853.4
$ws.Range("I107").Value = 4148.7144
$ws.Range("J107").Value = 3164.3333
$ws.Range("K107").Value = 12446.1432
$ws.Range("L107").Value = 9492.999899999999
$ws.Range("M107").Value = -10526.1432
$ws.Range("N107").Value = -13332.9999

# Row 132 - Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2595.125
$ws.Range("I132").Value = 2595.125
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7785.375
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5255.375
$ws.Range("N132").ClearContents()

# Row 134 - Cloth for Canvas / Mountain Linen
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
